$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 3 corresponds to the b889dd7b... file in each locale sheet.
$wsZhCn.Range("E3").Value = "2016-03-17 22:38:21"
$wsZhCn.Range("H3").Value = "2016-03-17 22:38:39"

$wsDeDe.Range("E3").Value = "2016-03-17 22:38:24"
$wsDeDe.Range("H3").Value = "2016-03-17 22:38:45"
